# chore: update Sheets via scheduled runner
#
# Applies the numeric refresh captured in the commit's xml diff across the
# ALC / ARM / BSM / CRP / CUL / GSM / LTW / WVR leve-profit tables: per-row
# currentAveragePrice* / LevePrice* / LeveProfit* columns (H:N) are
# recomputed, with a couple of cells cleared (no longer populated) or newly
# populated where the source row changed shape.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(19, 8).Value = 0
$ws.Cells.Item(19, 9).Value = 0
$ws.Cells.Item(19, 10).Value = 0
$ws.Cells.Item(19, 11).Value = 0
$ws.Cells.Item(19, 13).ClearContents()
$ws.Cells.Item(19, 14).ClearContents()
$ws.Cells.Item(100, 8).Value = 4192.593
$ws.Cells.Item(100, 10).Value = 5863.636
$ws.Cells.Item(100, 12).Value = 5863.636
$ws.Cells.Item(100, 14).Value = -6945.636
$ws.Cells.Item(113, 8).Value = 1452.3334
$ws.Cells.Item(113, 10).Value = 1466.875
$ws.Cells.Item(113, 12).Value = 1466.875
$ws.Cells.Item(113, 14).Value = -7974.875
$ws.Cells.Item(132, 8).Value = 151896.4
$ws.Cells.Item(132, 9).Value = 2394.3696
$ws.Cells.Item(132, 10).Value = 479377.06
$ws.Cells.Item(132, 11).Value = 7183.1088
$ws.Cells.Item(132, 12).Value = 1438131.18
$ws.Cells.Item(132, 13).Value = -4653.1088
$ws.Cells.Item(132, 14).Value = -1443191.18
$ws.Cells.Item(135, 8).Value = 8197552.5
$ws.Cells.Item(135, 9).Value = 251.8158
$ws.Cells.Item(135, 10).Value = 21740918
$ws.Cells.Item(135, 11).Value = 2266.3422
$ws.Cells.Item(135, 12).Value = 195668262
$ws.Cells.Item(135, 13).Value = 268.6578
$ws.Cells.Item(135, 14).Value = -195673332
$ws.Cells.Item(137, 8).Value = 40515.715
$ws.Cells.Item(137, 9).Value = 67934.8
$ws.Cells.Item(137, 10).Value = 8878.308000000001
$ws.Cells.Item(137, 11).Value = 203804.4
$ws.Cells.Item(137, 12).Value = 26634.924
$ws.Cells.Item(137, 13).Value = -201254.4
$ws.Cells.Item(137, 14).Value = -31734.924
$ws.Cells.Item(139, 8).Value = 50000
$ws.Cells.Item(139, 9).Value = 0
$ws.Cells.Item(139, 11).Value = 0
$ws.Cells.Item(139, 13).ClearContents()
$ws.Cells.Item(141, 8).Value = 1289.8928
$ws.Cells.Item(141, 9).Value = 716.5682
$ws.Cells.Item(141, 10).Value = 3392.0833
$ws.Cells.Item(141, 11).Value = 2149.7046
$ws.Cells.Item(141, 12).Value = 10176.2499
$ws.Cells.Item(141, 13).Value = 3030.2954
$ws.Cells.Item(141, 14).Value = -20536.2499
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 5466.2383
$ws.Cells.Item(2, 9).Value = 935.0625
$ws.Cells.Item(2, 10).Value = 19966
$ws.Cells.Item(2, 11).Value = 935.0625
$ws.Cells.Item(2, 12).Value = 19966
$ws.Cells.Item(2, 13).Value = -822.0625
$ws.Cells.Item(2, 14).Value = -20192
$ws.Cells.Item(11, 8).Value = 20000
$ws.Cells.Item(11, 10).Value = 20000
$ws.Cells.Item(11, 12).Value = 20000
$ws.Cells.Item(11, 14).Value = -20288
$ws.Cells.Item(32, 8).Value = 1595224.6
$ws.Cells.Item(32, 9).Value = 1969220.5
$ws.Cells.Item(32, 10).Value = 5742.0835
$ws.Cells.Item(32, 11).Value = 1969220.5
$ws.Cells.Item(32, 12).Value = 5742.0835
$ws.Cells.Item(32, 13).Value = -1968933.5
$ws.Cells.Item(32, 14).Value = -6316.0835
$ws.Cells.Item(74, 8).Value = 52513.3
$ws.Cells.Item(74, 9).Value = 103419.2
$ws.Cells.Item(74, 10).Value = 1607.4
$ws.Cells.Item(74, 11).Value = 103419.2
$ws.Cells.Item(74, 12).Value = 1607.4
$ws.Cells.Item(74, 13).Value = -102545.2
$ws.Cells.Item(74, 14).Value = -3355.4
$ws.Cells.Item(77, 8).Value = 52513.3
$ws.Cells.Item(77, 9).Value = 103419.2
$ws.Cells.Item(77, 10).Value = 1607.4
$ws.Cells.Item(77, 11).Value = 517096
$ws.Cells.Item(77, 12).Value = 8037
$ws.Cells.Item(77, 13).Value = -512728
$ws.Cells.Item(77, 14).Value = -16773
$ws.Cells.Item(116, 8).Value = 5466.2383
$ws.Cells.Item(116, 9).Value = 935.0625
$ws.Cells.Item(116, 10).Value = 19966
$ws.Cells.Item(116, 11).Value = 935.0625
$ws.Cells.Item(116, 12).Value = 19966
$ws.Cells.Item(116, 13).Value = 1358.9375
$ws.Cells.Item(116, 14).Value = -24554
$ws.Cells.Item(132, 8).Value = 2440482
$ws.Cells.Item(132, 9).Value = 3004218.2
$ws.Cells.Item(132, 10).Value = 843229.75
$ws.Cells.Item(132, 11).Value = 9012654.600000001
$ws.Cells.Item(132, 12).Value = 2529689.25
$ws.Cells.Item(132, 13).Value = -9010124.600000001
$ws.Cells.Item(132, 14).Value = -2534749.25
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 5466.2383
$ws.Cells.Item(3, 9).Value = 935.0625
$ws.Cells.Item(3, 10).Value = 19966
$ws.Cells.Item(3, 11).Value = 935.0625
$ws.Cells.Item(3, 12).Value = 19966
$ws.Cells.Item(3, 13).Value = -821.0625
$ws.Cells.Item(3, 14).Value = -20194
$ws.Cells.Item(105, 8).Value = 996301.25
$ws.Cells.Item(105, 9).Value = 1991127.5
$ws.Cells.Item(105, 10).Value = 1475
$ws.Cells.Item(105, 11).Value = 1991127.5
$ws.Cells.Item(105, 12).Value = 1475
$ws.Cells.Item(105, 13).Value = -1989380.5
$ws.Cells.Item(105, 14).Value = -4969
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 878.2
$ws.Cells.Item(16, 9).Value = 752.7143
$ws.Cells.Item(16, 10).Value = 1171
$ws.Cells.Item(16, 11).Value = 752.7143
$ws.Cells.Item(16, 12).Value = 1171
$ws.Cells.Item(16, 13).Value = -465.7143
$ws.Cells.Item(16, 14).Value = -1745
$ws.Cells.Item(31, 8).Value = 10959.069
$ws.Cells.Item(31, 9).Value = 9049.226000000001
$ws.Cells.Item(31, 10).Value = 15892.833
$ws.Cells.Item(31, 11).Value = 9049.226000000001
$ws.Cells.Item(31, 12).Value = 15892.833
$ws.Cells.Item(31, 13).Value = -8754.226000000001
$ws.Cells.Item(31, 14).Value = -16482.833
$ws.Cells.Item(34, 8).Value = 10959.069
$ws.Cells.Item(34, 9).Value = 9049.226000000001
$ws.Cells.Item(34, 10).Value = 15892.833
$ws.Cells.Item(34, 11).Value = 9049.226000000001
$ws.Cells.Item(34, 12).Value = 15892.833
$ws.Cells.Item(34, 13).Value = -8847.226000000001
$ws.Cells.Item(34, 14).Value = -16296.833
$ws.Cells.Item(94, 8).Value = 1828.96
$ws.Cells.Item(94, 9).Value = 2391.5
$ws.Cells.Item(94, 10).Value = 1453.9333
$ws.Cells.Item(94, 11).Value = 2391.5
$ws.Cells.Item(94, 12).Value = 1453.9333
$ws.Cells.Item(94, 13).Value = -1940.5
$ws.Cells.Item(94, 14).Value = -2355.9333
$ws.Cells.Item(99, 8).Value = 4654.4443
$ws.Cells.Item(99, 9).Value = 4340
$ws.Cells.Item(99, 10).Value = 5047.5
$ws.Cells.Item(99, 11).Value = 4340
$ws.Cells.Item(99, 12).Value = 5047.5
$ws.Cells.Item(99, 13).Value = -2842
$ws.Cells.Item(99, 14).Value = -8043.5
$ws.Cells.Item(113, 8).Value = 878.2
$ws.Cells.Item(113, 9).Value = 752.7143
$ws.Cells.Item(113, 10).Value = 1171
$ws.Cells.Item(113, 11).Value = 752.7143
$ws.Cells.Item(113, 12).Value = 1171
$ws.Cells.Item(113, 13).Value = 1417.2857
$ws.Cells.Item(113, 14).Value = -5511
$ws.Cells.Item(126, 8).Value = 4654.4443
$ws.Cells.Item(126, 9).Value = 4340
$ws.Cells.Item(126, 10).Value = 5047.5
$ws.Cells.Item(126, 11).Value = 13020
$ws.Cells.Item(126, 12).Value = 15142.5
$ws.Cells.Item(126, 13).Value = -10550
$ws.Cells.Item(126, 14).Value = -20082.5
$ws.Cells.Item(132, 8).Value = 1274.3334
$ws.Cells.Item(132, 9).Value = 817.7692
$ws.Cells.Item(132, 10).Value = 1899.1052
$ws.Cells.Item(132, 11).Value = 2453.3076
$ws.Cells.Item(132, 12).Value = 5697.3156
$ws.Cells.Item(132, 13).Value = 76.69239999999991
$ws.Cells.Item(132, 14).Value = -10757.3156
$ws.Cells.Item(134, 8).Value = 1095.386
$ws.Cells.Item(134, 9).Value = 988.6818
$ws.Cells.Item(134, 10).Value = 1456.5385
$ws.Cells.Item(134, 11).Value = 2966.0454
$ws.Cells.Item(134, 12).Value = 4369.6155
$ws.Cells.Item(134, 13).Value = -431.0454
$ws.Cells.Item(134, 14).Value = -9439.6155
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(102, 8).Value = 2581.1428
$ws.Cells.Item(102, 10).Value = 2581.1428
$ws.Cells.Item(102, 12).Value = 7743.428400000001
$ws.Cells.Item(102, 14).Value = -12611.4284
$ws.Cells.Item(113, 8).Value = 521
$ws.Cells.Item(113, 9).Value = 450
$ws.Cells.Item(113, 10).Value = 556.5
$ws.Cells.Item(113, 11).Value = 1350
$ws.Cells.Item(113, 12).Value = 1669.5
$ws.Cells.Item(113, 13).Value = 820
$ws.Cells.Item(113, 14).Value = -6009.5
$ws.Cells.Item(120, 8).Value = 12766.294
$ws.Cells.Item(120, 9).Value = 5507.5
$ws.Cells.Item(120, 10).Value = 14999.77
$ws.Cells.Item(120, 11).Value = 16522.5
$ws.Cells.Item(120, 12).Value = 44999.31
$ws.Cells.Item(120, 13).Value = -11684.5
$ws.Cells.Item(120, 14).Value = -54675.31
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(107, 8).Value = 11166.889
$ws.Cells.Item(107, 9).Value = 312.75
$ws.Cells.Item(107, 10).Value = 98000
$ws.Cells.Item(107, 11).Value = 312.75
$ws.Cells.Item(107, 12).Value = 98000
$ws.Cells.Item(107, 13).Value = 1607.25
$ws.Cells.Item(107, 14).Value = -101840
$ws.Cells.Item(126, 8).Value = 1326.1904
$ws.Cells.Item(126, 9).Value = 1290.909
$ws.Cells.Item(126, 10).Value = 1455.5555
$ws.Cells.Item(126, 11).Value = 3872.727
$ws.Cells.Item(126, 12).Value = 4366.666499999999
$ws.Cells.Item(126, 13).Value = -1402.727
$ws.Cells.Item(126, 14).Value = -9306.666499999999
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(2, 8).Value = 50001.332
$ws.Cells.Item(2, 10).Value = 50001.332
$ws.Cells.Item(2, 12).Value = 50001.332
$ws.Cells.Item(2, 14).Value = -50225.332
$ws.Cells.Item(132, 8).Value = 225365.78
$ws.Cells.Item(132, 9).Value = 57711.39
$ws.Cells.Item(132, 10).Value = 560674.5600000001
$ws.Cells.Item(132, 11).Value = 173134.17
$ws.Cells.Item(132, 12).Value = 1682023.68
$ws.Cells.Item(132, 13).Value = -170604.17
$ws.Cells.Item(132, 14).Value = -1687083.68
$ws.Cells.Item(136, 8).Value = 456841.28
$ws.Cells.Item(136, 9).Value = 716077.7
$ws.Cells.Item(136, 10).Value = 3177.5
$ws.Cells.Item(136, 11).Value = 2148233.1
$ws.Cells.Item(136, 12).Value = 9532.5
$ws.Cells.Item(136, 13).Value = -2145683.1
$ws.Cells.Item(136, 14).Value = -14632.5
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(81, 8).Value = 983.6799999999999
$ws.Cells.Item(81, 9).Value = 652.7692
$ws.Cells.Item(81, 10).Value = 1342.1666
$ws.Cells.Item(81, 11).Value = 1305.5384
$ws.Cells.Item(81, 12).Value = 2684.3332
$ws.Cells.Item(81, 13).Value = -244.5383999999999
$ws.Cells.Item(81, 14).Value = -4806.3332
$ws.Cells.Item(84, 8).Value = 983.6799999999999
$ws.Cells.Item(84, 9).Value = 652.7692
$ws.Cells.Item(84, 10).Value = 1342.1666
$ws.Cells.Item(84, 11).Value = 6527.691999999999
$ws.Cells.Item(84, 12).Value = 13421.666
$ws.Cells.Item(84, 13).Value = -1223.691999999999
$ws.Cells.Item(84, 14).Value = -24029.666
$ws.Cells.Item(107, 8).Value = 359
$ws.Cells.Item(107, 9).Value = 373
$ws.Cells.Item(107, 10).Value = 303
$ws.Cells.Item(107, 11).Value = 1119
$ws.Cells.Item(107, 12).Value = 909
$ws.Cells.Item(107, 13).Value = 801
$ws.Cells.Item(107, 14).Value = -4749
$ws.Cells.Item(132, 8).Value = 3033.7292
$ws.Cells.Item(132, 9).Value = 640.70966
$ws.Cells.Item(132, 10).Value = 7397.4707
$ws.Cells.Item(132, 11).Value = 1922.12898
$ws.Cells.Item(132, 12).Value = 22192.4121
$ws.Cells.Item(132, 13).Value = 607.87102
$ws.Cells.Item(132, 14).Value = -27252.4121
